# useraccount.jsp now displays all orders per user
#
# 1. Change the highlight on the "Page listing all orders for user" row
#    (Feature table, row describing that feature) from yellow to
#    bright green.
# 2. Mark that feature as implemented by putting "1" in the "Points"
#    cell that was previously empty.

$d = $word.ActiveDocument

# Locate the table row for the "Page listing all orders for user" feature.
$targetTable = $null
$targetRow = 0
for ($t = 1; $t -le $d.Tables.Count; $t++) {
    $tbl = $d.Tables.Item($t)
    for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
        $cellText = $tbl.Cell($r, 1).Range.Text
        if ($cellText -like "*Page listing all orders for user*") {
            $targetTable = $tbl
            $targetRow = $r
            break
        }
    }
    if ($targetTable -ne $null) {
        break
    }
}

# 1) Re-highlight the feature-name cell in bright green (wdBrightGreen = 4).
$nameCell = $targetTable.Cell($targetRow, 1)
$nameCell.Range.Font.HighlightColorIndex = 4

# 2) Fill in the "Points" cell (column 4) with "1".
$pointsCell = $targetTable.Cell($targetRow, 4)
$pointsCell.Range.Text = "1"
